$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on the existing sheets -------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the same header/date formatting already used on the Weekly sheet so
# the new sheet's styles line up with the workbook's existing style table.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Header row -----------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 4. Data rows --------------------------------------------------------
$wsForecast.Range("A2").Value = 45123.99999999999
$wsForecast.Range("B2").Value = 1
$wsForecast.Range("C2").Value = 1.000009447400555
$wsForecast.Range("D2").Value = 1.000009447725985

$wsForecast.Range("A3").Value = 45179.99999999999
$wsForecast.Range("B3").Value = 10
$wsForecast.Range("C3").Value = 10.00000944943508
$wsForecast.Range("D3").Value = 10.00000944977175

$wsForecast.Range("A4").Value = 45186.99999999999
$wsForecast.Range("B4").Value = 11
$wsForecast.Range("C4").Value = 11.12500944965423
$wsForecast.Range("D4").Value = 11.12500945007576

$wsForecast.Range("A5").Value = 45193.99999999999
$wsForecast.Range("B5").Value = 12
$wsForecast.Range("C5").Value = 12.25000944978065
$wsForecast.Range("D5").Value = 12.25000945353299

$wsForecast.Range("A6").Value = 45200.99999999999
$wsForecast.Range("B6").Value = 13
$wsForecast.Range("C6").Value = 13.37500944524445
$wsForecast.Range("D6").Value = 13.37500946122255

$wsForecast.Range("A7").Value = 45207.99999999999
$wsForecast.Range("B7").Value = 15
$wsForecast.Range("C7").Value = 14.50000943743317
$wsForecast.Range("D7").Value = 14.50000947048582

$wsForecast.Range("A8").Value = 45214.99999999999
$wsForecast.Range("B8").Value = 16
$wsForecast.Range("C8").Value = 15.62500942657035
$wsForecast.Range("D8").Value = 15.62500948468458

$wsForecast.Range("A9").Value = 45221.99999999999
$wsForecast.Range("B9").Value = 17
$wsForecast.Range("C9").Value = 16.75000941473392
$wsForecast.Range("D9").Value = 16.75000949917319

$wsForecast.Range("A10").Value = 45228.99999999999
$wsForecast.Range("B10").Value = 18
$wsForecast.Range("C10").Value = 17.87500939970966
$wsForecast.Range("D10").Value = 17.87500951584322

$wsForecast.Range("A11").Value = 45235.99999999999
$wsForecast.Range("B11").Value = 19
$wsForecast.Range("C11").Value = 19.00000937875809
$wsForecast.Range("D11").Value = 19.00000953223734

# Put the selection back on the first sheet, matching the original file.
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select()
